# edit.ps1 - apply the two text corrections described by the diff:
#  1) Slide 3, "Content Placeholder 2": merge/fix
#       "Those are objects provided by the Google VR SDK to take care of core VR "
#       + "functionnalities"
#     into a single corrected run:
#       "Those are objects provided by the Google VR SDK to take care of core VR functionalities"
#  2) Slide 5, "Content Placeholder 2": fix
#       "However added functionalities are added to the camera"
#     into:
#       "However additional functionalities are added to the camera"
#
# Note: target slide indices are passed as hints but every slide is searched
# as a fallback, so the script is resilient to slides being renumbered.

function Fix-ShapeText {
    param($Shape, $OldText, $NewText)

    if (-not $Shape.HasTextFrame) { return $false }
    if (-not $Shape.TextFrame.HasText) { return $false }

    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        return $false
    }

    $sub = $tr.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
    return $true
}

function Fix-SlideText {
    param($Slide, $OldText, $NewText)

    $done = $false
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shape = $Slide.Shapes.Item($i)
        if (Fix-ShapeText $shape $OldText $NewText) {
            $done = $true
            break
        }
    }
    return $done
}

function Fix-PresentationText {
    param($Presentation, $PreferredSlideIndex, $OldText, $NewText)

    # Try the preferred slide first (fast path / matches the original deck).
    if ($PreferredSlideIndex -ge 1 -and $PreferredSlideIndex -le $Presentation.Slides.Count) {
        $slide = $Presentation.Slides.Item($PreferredSlideIndex)
        if (Fix-SlideText $slide $OldText $NewText) {
            return $true
        }
    }

    # Fallback: scan every slide.
    for ($s = 1; $s -le $Presentation.Slides.Count; $s++) {
        $slide = $Presentation.Slides.Item($s)
        if (Fix-SlideText $slide $OldText $NewText) {
            return $true
        }
    }

    return $false
}

$p = $ppt.ActivePresentation

# --- fix the "core VR functionnalities" typo (also merges the two runs that
#     made up the sentence into the single corrected run) -------------------
$old3 = "Those are objects provided by the Google VR SDK to take care of core VR functionnalities"
$new3 = "Those are objects provided by the Google VR SDK to take care of core VR functionalities"
Fix-PresentationText $p 3 $old3 $new3

# --- fix "However added functionalities" -> "However additional functionalities"
$old5 = "However added functionalities are added to the camera"
$new5 = "However additional functionalities are added to the camera"
Fix-PresentationText $p 5 $old5 $new5
